$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.824999999999998
$ws.Range("D8").Value = -8.048999999999999
$ws.Range("B12").Value = 6.112
$ws.Range("D12").Value = -8.311
$ws.Range("D14").Value = -8.263
$ws.Range("D22").Value = -8.191999999999998
